$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B, shifting the existing
# Jun_17 / Jun_15 / Jun_13 / Jun_10 columns (old B:E) to E:H.
$ws.Columns("B:D").Insert()

# Give columns C-H (raw OOXML width 8.0 == ColumnWidth 7.1667) a
# matching custom width; column B keeps the default width.
$ws.Columns("C").ColumnWidth = 7.1667
$ws.Columns("D").ColumnWidth = 7.1667
$ws.Columns("E").ColumnWidth = 7.1667
$ws.Columns("F").ColumnWidth = 7.1667
$ws.Columns("G").ColumnWidth = 7.1667
$ws.Columns("H").ColumnWidth = 7.1667

# New header row values for the newly inserted date columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the new columns for every existing data row (2-27) with the
# same "UN" placeholder used throughout the rest of the table.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# New group rows appended at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
